$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated TPM-derived values (recomputed by upstream NATMI script) for rows 2-10.
# Each assignment below mirrors a single cell change from the commit diff.
$ws.Range("G2").Value = 21.65036466666666
$ws.Range("H2").Value = 64.951094
$ws.Range("I2").Value = 0.9284333993050746
$ws.Range("J2").Value = 0.9284333993050747
$ws.Range("M2").Value = 72.07569866666667
$ws.Range("N2").Value = 216.227096
$ws.Range("O2").Value = 0.4479522040449755
$ws.Range("P2").Value = 0.4479522040449755
$ws.Range("Q2").Value = 1560.465159738114
$ws.Range("R2").Value = 14044.18643764302
$ws.Range("S2").Value = 0.415893787527677
$ws.Range("T2").Value = 0.4158937875276771
$ws.Range("G3").Value = 21.65036466666666
$ws.Range("H3").Value = 64.951094
$ws.Range("I3").Value = 0.9284333993050746
$ws.Range("J3").Value = 0.9284333993050747
$ws.Range("O3").Value = 0.04737448730867841
$ws.Range("P3").Value = 0.0473744873086784
$ws.Range("Q3").Value = 165.0315284490167
$ws.Range("R3").Value = 1485.28375604115
$ws.Range("S3").Value = 0.04398405629233141
$ws.Range("T3").Value = 0.04398405629233141
$ws.Range("G4").Value = 21.65036466666666
$ws.Range("H4").Value = 64.951094
$ws.Range("I4").Value = 0.9284333993050746
$ws.Range("J4").Value = 0.9284333993050747
$ws.Range("O4").Value = 0.5046733086463462
$ws.Range("P4").Value = 0.5046733086463461
$ws.Range("Q4").Value = 1758.056123133425
$ws.Range("R4").Value = 15822.50510820083
$ws.Range("S4").Value = 0.4685555554850663
$ws.Range("T4").Value = 0.4685555554850663
$ws.Range("G5").Value = 0.96805
$ws.Range("I5").Value = 0.04151292442575075
$ws.Range("J5").Value = 0.04151292442575075
$ws.Range("M5").Value = 72.07569866666667
$ws.Range("N5").Value = 216.227096
$ws.Range("O5").Value = 0.4479522040449755
$ws.Range("P5").Value = 0.4479522040449755
$ws.Range("Q5").Value = 69.77288009426667
$ws.Range("R5").Value = 627.9559208484001
$ws.Range("S5").Value = 0.01859580599286755
$ws.Range("T5").Value = 0.01859580599286755
$ws.Range("G6").Value = 0.96805
$ws.Range("I6").Value = 0.04151292442575075
$ws.Range("J6").Value = 0.04151292442575075
$ws.Range("O6").Value = 0.04737448730867841
$ws.Range("P6").Value = 0.0473744873086784
$ws.Range("Q6").Value = 7.37903372875
$ws.Range("S6").Value = 0.001966653511353855
$ws.Range("T6").Value = 0.001966653511353854
$ws.Range("G7").Value = 0.96805
$ws.Range("I7").Value = 0.04151292442575075
$ws.Range("J7").Value = 0.04151292442575075
$ws.Range("O7").Value = 0.5046733086463462
$ws.Range("P7").Value = 0.5046733086463461
$ws.Range("Q7").Value = 78.60773969408334
$ws.Range("R7").Value = 707.46965724675
$ws.Range("S7").Value = 0.02095046492152935
$ws.Range("T7").Value = 0.02095046492152934
$ws.Range("I8").Value = 0.0300536762691746
$ws.Range("J8").Value = 0.0300536762691746
$ws.Range("M8").Value = 72.07569866666667
$ws.Range("N8").Value = 216.227096
$ws.Range("O8").Value = 0.4479522040449755
$ws.Range("P8").Value = 0.4479522040449755
$ws.Range("Q8").Value = 50.51273982086133
$ws.Range("R8").Value = 454.6146583877521
$ws.Range("S8").Value = 0.01346261052443094
$ws.Range("T8").Value = 0.01346261052443094
$ws.Range("I9").Value = 0.0300536762691746
$ws.Range("J9").Value = 0.0300536762691746
$ws.Range("O9").Value = 0.04737448730867841
$ws.Range("P9").Value = 0.0473744873086784
$ws.Range("S9").Value = 0.001423777504993142
$ws.Range("T9").Value = 0.001423777504993141
$ws.Range("I10").Value = 0.0300536762691746
$ws.Range("J10").Value = 0.0300536762691746
$ws.Range("O10").Value = 0.5046733086463462
$ws.Range("P10").Value = 0.5046733086463461
$ws.Range("S10").Value = 0.01516728823975052
$ws.Range("T10").Value = 0.01516728823975052
